# Auto-generated Excel COM-interop script
# Applies the cryptos.xlsx price/volume update described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new value is a numeric-looking string (e.g. "1.00", "7.83") ---
# These must be pre-formatted as Text so Excel stores the exact literal string
# instead of silently converting it to a number (which would drop trailing zeros,
# e.g. "1.00" -> 1).
$numericLookingTextCells = [ordered]@{
    "D5" = "604.69"
    "D6" = "144.39"
    "D8" = "1.00"
    "D10" = "7.83"
    "D12" = "0.409"
    "D15" = "28.56"
    "D19" = "11.10"
    "D20" = "6.17"
    "D21" = "14.62"
    "D22" = "422.08"
    "D24" = "77.21"
    "D28" = "7.89"
    "D29" = "2.46"
    "D30" = "8.92"
    "D31" = "1.00"
    "D34" = "24.24"
    "D37" = "7.61"
    "D38" = "1.65"
    "D39" = "175.63"
    "D40" = "5.25"
    "D41" = "0.0821"
    "D42" = "0.861"
    "D43" = "4.99"
    "D44" = "45.45"
    "D48" = "23.77"
    "D49" = "7.06"
    "D50" = "1.12"
    "D51" = "0.911"
}
foreach ($ref in $numericLookingTextCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingTextCells[$ref]
}

# --- Remaining cells: plain text / non-numeric-looking values ---
$plainCells = [ordered]@{
    "D2" = "66.365.40"
    "E2" = "  -0.29%  "
    "D3" = "3.538.10"
    "E3" = "  -1.44%  "
    "E4" = "  +0.01%  "
    "E5" = "  -0.71%  "
    "E6" = "  -2.47%  "
    "D7" = "3.538.60"
    "E7" = "  -1.38%  "
    "E8" = "  -0.11%  "
    "E9" = "  +5.15%  "
    "E10" = "  -2.60%  "
    "E11" = "  -3.76%  "
    "E12" = "  -1.36%  "
    "D13" = "4.145.23"
    "E13" = "  -1.36%  "
    "E14" = "  -6.68%  "
    "E15" = "  -4.77%  "
    "D16" = "3.540.36"
    "E16" = "  -1.62%  "
    "E17" = "  +0.98%  "
    "D18" = "66.259.16"
    "E19" = "  -3.54%  "
    "E20" = "  -2.92%  "
    "E21" = "  -2.89%  "
    "E22" = "  -1.49%  "
    "E23" = "  -3.94%  "
    "E24" = "  -2.31%  "
    "D25" = "3.679.32"
    "E25" = "  -1.60%  "
    "E26" = "  +0.06%  "
    "E27" = "  -5.18%  "
    "E28" = "  -4.63%  "
    "E29" = "  -2.13%  "
    "E30" = "  -3.97%  "
    "E31" = "  -0.02%  "
    "D32" = "3.546.60"
    "E32" = "  -1.18%  "
    "E33" = "  -0.91%  "
    "E34" = "  -4.83%  "
    "E35" = "  -0.04%  "
    "E36" = "  -7.45%  "
    "E37" = "  -2.89%  "
    "E38" = "  -4.23%  "
    "E39" = "  -0.65%  "
    "E40" = "  -6.83%  "
    "E41" = "  -4.38%  "
    "B42" = "Mantle"
    "C42" = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
    "E42" = "  -4.02%  "
    "B43" = "Filecoin"
    "C43" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "E43" = "  -4.65%  "
    "E44" = "  -1.61%  "
    "E45" = "  -7.78%  "
    "E46" = "  +0.14%  "
    "E47" = "  -6.72%  "
    "E48" = "  -1.88%  "
    "E49" = "  -1.74%  "
    "E50" = "  -6.09%  "
    "E51" = "  -4.19%  "
}
foreach ($ref in $plainCells.Keys) {
    $ws.Range($ref).Value = $plainCells[$ref]
}

